$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7743191123008728
$ws.Range("B1").Value = 1.185871481895447
$ws.Range("C1").Value = 2.314530372619629
$ws.Range("D1").Value = -1
$ws.Range("E1").Value = 1.794419169425964
